# Generate Report for Handback
# Update the "generate/handoff/handback" timestamp cells to reflect a new
# report-generation run. All of these columns hold their datetime as plain
# text (they are shared-string cells formatted with a date display mask),
# so we force the assignment to text to avoid Excel re-interpreting it as
# a date serial number.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$overview.Range("G2").Value = "2016-08-18 23:06:12"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2)
$zhcn.Range("H2").Value = "2016-08-18 23:06:05"
$zhcn.Range("K2").Value = "2016-08-18 23:06:31"

# de-de sheet: "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2)
$dede.Range("H2").Value = "2016-08-18 23:06:12"
$dede.Range("K2").Value = "2016-08-18 23:06:39"
